$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old range (A2:A6) entirely
$ws.Range("A2:A6").Clear()

# Write the new values into A1:A5
$ws.Range("A1").Value = "LRfY"
$ws.Range("A2").Value = "qXC Jc"
$ws.Range("A3").Value = "NpLW0dH5Q0EFkF"
$ws.Range("A4").Value = "f lmRU"
$ws.Range("A5").Value = "RlGpDG4 WYj"
